$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.423.99"

# Row 3
$ws.Range("D3").Value = "1.945.68"
$ws.Range("E3").Value = "  -1.85%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'242.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "

# Row 6
$ws.Range("D6").Value = "'0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.28%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").Value = "'56.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.51%  "

# Row 9
$ws.Range("D9").Value = "'0.361"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.97%  "

# Row 10
$ws.Range("E10").Value = "  +3.16%  "

# Row 11
$ws.Range("E11").Value = "  +0.01%  "

# Row 12
$ws.Range("D12").Value = "2.233.14"
$ws.Range("E12").Value = "  -1.68%  "

# Row 13
$ws.Range("E13").Value = "  -5.70%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'13.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.88%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'21.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -11.58%  "

# Row 16
$ws.Range("D16").Value = "'5.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.50%  "

# Row 17
$ws.Range("D17").Value = "1.957.15"
$ws.Range("E17").Value = "  -1.27%  "

# Row 18
$ws.Range("D18").Value = "36.336.15"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0872"
$ws.Range("E19").Value = "  +1.23%  "

# Row 20
$ws.Range("D20").Value = "'69.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.80%  "

# Row 21
$ws.Range("D21").Value = "'228.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.41%  "

# Row 22
$ws.Range("D22").Value = "'5.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.86%  "

# Row 23
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$ws.Range("E24").Value = "  -8.01%  "

# Row 25
$ws.Range("E25").Value = "  -1.09%  "

# Row 26
$ws.Range("D26").Value = "'9.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.48%  "

# Row 27
$ws.Range("D27").Value = "'161.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.21%  "

# Row 28
$ws.Range("E28").Value = "  +5.92%  "

# Row 29
$ws.Range("D29").Value = "'19.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.29%  "

# Row 30
$ws.Range("D30").Value = "'0.117"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.22%  "

# Row 31
$ws.Range("E31").Value = "  -4.83%  "

# Row 32
$ws.Range("D32").Value = "'4.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.14%  "

# Row 33
$ws.Range("D33").Value = "'0.0635"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.55%  "

# Row 34
$ws.Range("E34").Value = "  -3.82%  "

# Row 35
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "

# Row 36
$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").Value = "'6.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.04%  "

# Row 37
$ws.Range("E37").Value = "  +1.16%  "

# Row 38
$ws.Range("D38").Value = "'2.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.03%  "

# Row 39
$ws.Range("D39").Value = "'3.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.46%  "

# Row 40
$ws.Range("D40").Value = "'0.0967"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "

# Row 41
$ws.Range("D41").Value = "'2.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.88%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0210"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.69%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.12%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'15.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.99%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.346.63"
$ws.Range("E45").Value = "  -2.22%  "

# Row 46
$ws.Range("E46").Value = "  -6.62%  "

# Row 47
$ws.Range("D47").Value = "'87.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.97%  "

# Row 48
$ws.Range("D48").Value = "'7.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.23%  "

# Row 49
$ws.Range("E49").Value = "  -0.73%  "

# Row 50
$ws.Range("D50").Value = "'44.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.79%  "

# Row 51
$ws.Range("D51").Value = "2.123.75"
$ws.Range("E51").Value = "  -1.92%  "
